$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Design Parameters"

# --- Capture existing threaded-comment text for C1/D1/E1 before shifting them ---
$txtMin = $ws.Range("C1").CommentThreaded.Text()
$txtMax = $ws.Range("D1").CommentThreaded.Text()
$txtDrop = $ws.Range("E1").CommentThreaded.Text()

# Remove the old threaded comments (also clears the legacy mirror comment)
$ws.Range("C1").CommentThreaded.Delete()
$ws.Range("D1").CommentThreaded.Delete()
$ws.Range("E1").CommentThreaded.Delete()

# --- Insert the new "Feature type" column before the old column C ---
$ws.Columns.Item(3).Insert()

# New header + data for the inserted column
$ws.Range("C1").Value = "Feature type"
$ws.Range("C2").Value = "Numerical"
$ws.Range("C3").Value = "Numerical"
$ws.Range("C4").Value = "Numerical"

# Match formatting used by its neighbours
$ws.Range("C1").Style = $ws.Range("D1").Style
$ws.Range("C2:C4").Style = $ws.Range("B2:B4").Style

# Column width close to the authored width
$ws.Columns.Item(3).ColumnWidth = 11.8

# Data validation dropdown for the new column
$ws.Range("C2:C4").Validation.Add(3, 1, 1, '"Numerical, Categorical"')

# Re-add the threaded comments at their shifted locations
$ws.Range("D1").AddCommentThreaded($txtMin)
$ws.Range("E1").AddCommentThreaded($txtMax)
$ws.Range("F1").AddCommentThreaded($txtDrop)

# Selection / active sheet bookkeeping (matches the authored file)
$ws.Range("C3:C4").Select()
$ws.Activate()

# --- "Train Data" sheet loses the tabSelected flag when Design Parameters becomes active ---
# (handled automatically by activating the Design Parameters sheet above)
